$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a "last updated" date serial number.
# It is bumped by one day (45189 -> 45190, i.e. 2023-09-20 -> 2023-09-21)
# for every populated data row (rows 2 through 98).
for ($r = 2; $r -le 98; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
